# Updated Data base verification test cases
$wb = $excel.ActiveWorkbook

# --- validLoginTest (sheet1): refresh the credentials used by the test ---
$ws1 = $wb.Worksheets.Item("validLoginTest")
$ws1.Range("A2").Value = "orangehrm_vani"
$ws1.Range("B2").Value = "Qweinav12!8"
$ws1.Range("B2").Select() | Out-Null

# --- add the new empDBVerification sheet as the last (3rd) tab and activate it ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "empDBVerification"

$ws3.Range("A1").Value = "empID"
$ws3.Range("B1").Value = "empName"

$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "Vani Bhat"

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "ashwin hebbar"

$ws3.Range("A4").Value = 3
$ws3.Range("B4").Value = "sasha de della"

$ws3.Columns.Item(1).ColumnWidth = 11
$ws3.Columns.Item(2).ColumnWidth = 12.333333333333334

$ws3.Range("B4").Select() | Out-Null
